# "Set up camera system" - append the February 26th devlog entry.
#
# Target shape (mirrors the existing "February 25th:" entry a few
# paragraphs above it):
#   <empty paragraph>
#   February 26th:                       (bold heading, "th" superscript)
#   Camera system set up, can swap between states. Current states
#   include main and keg, but more can be added as needed.
#   Cup will need to be dragged to location under tap?

$d = $word.ActiveDocument

# Collapse to the very end of the document body, right after
# "Keg now checks if tray contains cup when clicked."
$end = $d.Content
$end.Collapse(0)

# Build the new paragraphs as a WordProcessingML fragment (the same
# "single file package" shape Range.WordOpenXML round-trips) so the
# empty spacer paragraph comes through with no stray run, exactly like
# the pre-existing blank paragraph before "February 25th:".
$newParagraphsXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>February 26</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Camera system set up, can swap between states. Current states include main and keg, but more can be added as needed.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Cup will need to be dragged to location under tap?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$end.InsertXML($newParagraphsXml)
